# Natmi following Dr Hou advice
# Rebuild the LR-pair table: re-derive per-sender-cluster aggregate stats
# (columns E:T) now that "sCs" and "FAPs" are both treated as separate
# sending clusters, and append the new FAPs-as-sender rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E:T = the 16 numeric LR-pair statistics.
$data = @(
    @("FAPs", "Il12a", "Il12rb2", "ECs", 1, 0.3333333333333333, 0.3206746666666667, 0.962024, 0.3143489699528113, 0.4074782680148976, 3, 1, 0.2503153333333333, 0.750946, 0.0905543021210042, 0.09435844210674453, 0.08026978607822222, 0.7224280747039999, 0.02846565159653334, 0.03844901456224025),
    @("FAPs", "Il12a", "Il12rb2", "FAPs", 1, 0.3333333333333333, 0.3206746666666667, 0.962024, 0.3143489699528113, 0.4074782680148976, 3, 1, 0.21436, 0.64308, 0.07754706810872602, 0.08080478083644532, 0.06873982154666666, 0.61865839392, 0.02437684098283853, 0.03292619214255813),
    @("FAPs", "Il12a", "Il12rb2", "M1", 1, 0.3333333333333333, 0.3206746666666667, 0.962024, 0.3143489699528113, 0.4074782680148976, 3, 1, 0.8543153333333334, 2.562946, 0.3090578901862707, 0.3220412543161725, 0.2739572847448889, 2.465615562704, 0.09715202943584325, 0.1312248125380991),
    @("FAPs", "Il12a", "Il12rb2", "M2", 1, 0.3333333333333333, 0.3206746666666667, 0.962024, 0.3143489699528113, 0.4074782680148976, 3, 1, 1.110936, 3.332808, 0.4018932154153557, 0.4187765441468428, 0.356249031488, 3.206241283392, 0.1263347182968403, 0.1706423408942198),
    @("FAPs", "Il12a", "Il12rb2", "sCs", 1, 0.3333333333333333, 0.3206746666666667, 0.962024, 0.3143489699528113, 0.4074782680148976, 2, 1, 0.33433, 0.66866, 0.1209475241686433, 0.08401897859379476, 0.1072111613066667, 0.64326696784, 0.03801972964075576, 0.03423590787778025),
    @("sCs", "Il12a", "Il12rb2", "ECs", 2, 1, 0.6994484999999999, 1.398897, 0.6856510300471887, 0.5925217319851024, 3, 1, 0.2503153333333333, 0.750946, 0.0905543021210042, 0.09435844210674453, 0.175082684427, 1.050496106562, 0.06208865052447086, 0.05590942754450428),
    @("sCs", "Il12a", "Il12rb2", "FAPs", 2, 1, 0.6994484999999999, 1.398897, 0.6856510300471887, 0.5925217319851024, 3, 1, 0.21436, 0.64308, 0.07754706810872602, 0.08080478083644532, 0.14993378046, 0.8996026827599999, 0.0531702271258875, 0.04787858869388719),
    @("sCs", "Il12a", "Il12rb2", "M1", 2, 1, 0.6994484999999999, 1.398897, 0.6856510300471887, 0.5925217319851024, 3, 1, 0.8543153333333334, 2.562946, 0.3090578901862707, 0.3220412543161725, 0.597549578427, 3.585297470562, 0.2119058607504275, 0.1908164417780734),
    @("sCs", "Il12a", "Il12rb2", "M2", 2, 1, 0.6994484999999999, 1.398897, 0.6856510300471887, 0.5925217319851024, 3, 1, 1.110936, 3.332808, 0.4018932154153557, 0.4187765441468428, 0.7770425187959998, 4.662255112775999, 0.2755584971185154, 0.248134203252623),
    @("sCs", "Il12a", "Il12rb2", "sCs", 2, 1, 0.6994484999999999, 1.398897, 0.6856510300471887, 0.5925217319851024, 2, 1, 0.33433, 0.66866, 0.1209475241686433, 0.08401897859379476, 0.233846617005, 0.9353864680199999, 0.08292779452788752, 0.04978307071601451)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $row[$j]
    }
}
